# Generate Report for Handoff
#
# Swaps the "26532f2a-759a-4e9d-abbc-6e4fefee44ae.md" and
# "502c27b6-88d5-46b1-aef4-3ab37a07d94c.md" data rows (502c27b6 now listed
# first / row 2, 26532f2a now listed second / row 3) across the Overview,
# zh-cn and de-de sheets, and marks 26532f2a as ready for handoff with a
# fresh "Latest Handoff Datetime" in each locale sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "502c27b6-88d5-46b1-aef4-3ab37a07d94c.md"
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "In Translation"

$ws.Range("A3").Value = "26532f2a-759a-4e9d-abbc-6e4fefee44ae.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "502c27b6-88d5-46b1-aef4-3ab37a07d94c.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "26532f2a-759a-4e9d-abbc-6e4fefee44ae.md"
    }
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "502c27b6-88d5-46b1-aef4-3ab37a07d94c.md"
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "502c27b6-88d5-46b1-aef4-3ab37a07d94c.710498b3a79c43a0bd46a7918448d91bd3ed30ac.zh-cn.xlf"
$ws.Range("D2").Value = "2016-03-09 00:14:25"
$ws.Range("G2").Value = "0001-01-01 00:00:00"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "26532f2a-759a-4e9d-abbc-6e4fefee44ae.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "26532f2a-759a-4e9d-abbc-6e4fefee44ae.92ee07ea2c7bf5d31c35e785427dceee86103759.zh-cn.xlf"
$ws.Range("D3").Value = "2016-03-09 00:16:36"
$ws.Range("G3").Value = "0001-01-01 00:00:00"
$ws.Range("H3").Value = "Include"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "502c27b6-88d5-46b1-aef4-3ab37a07d94c.md"
    } elseif ($addr -eq '$C$2') {
        $hl.TextToDisplay = "502c27b6-88d5-46b1-aef4-3ab37a07d94c.710498b3a79c43a0bd46a7918448d91bd3ed30ac.zh-cn.xlf"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "26532f2a-759a-4e9d-abbc-6e4fefee44ae.md"
    } elseif ($addr -eq '$C$3') {
        $hl.TextToDisplay = "26532f2a-759a-4e9d-abbc-6e4fefee44ae.92ee07ea2c7bf5d31c35e785427dceee86103759.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "502c27b6-88d5-46b1-aef4-3ab37a07d94c.md"
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "502c27b6-88d5-46b1-aef4-3ab37a07d94c.710498b3a79c43a0bd46a7918448d91bd3ed30ac.de-de.xlf"
$ws.Range("D2").Value = "2016-03-09 00:15:08"
$ws.Range("G2").Value = "0001-01-01 00:00:00"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "26532f2a-759a-4e9d-abbc-6e4fefee44ae.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "26532f2a-759a-4e9d-abbc-6e4fefee44ae.92ee07ea2c7bf5d31c35e785427dceee86103759.de-de.xlf"
$ws.Range("D3").Value = "2016-03-09 00:16:45"
$ws.Range("G3").Value = "0001-01-01 00:00:00"
$ws.Range("H3").Value = "Include"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "502c27b6-88d5-46b1-aef4-3ab37a07d94c.md"
    } elseif ($addr -eq '$C$2') {
        $hl.TextToDisplay = "502c27b6-88d5-46b1-aef4-3ab37a07d94c.710498b3a79c43a0bd46a7918448d91bd3ed30ac.de-de.xlf"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "26532f2a-759a-4e9d-abbc-6e4fefee44ae.md"
    } elseif ($addr -eq '$C$3') {
        $hl.TextToDisplay = "26532f2a-759a-4e9d-abbc-6e4fefee44ae.92ee07ea2c7bf5d31c35e785427dceee86103759.de-de.xlf"
    }
}
